# Update the SMOTE model comparison table: rows 3 and 4 (B..E) get new
# metric values. Columns B and E are numeric; columns C and D are stored
# as text in the source sheet, so we force text storage (NumberFormat
# "@") before assigning the string value, then restore the default
# "Normal" style so no stray formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("Financial Variables and Sector")
$ws.Range("B3").Value = 0.9472
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "0.9097"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0.9472"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = 0.9281

# Row 4 ("Financial Variables, Sector, and NLP Features")
$ws.Range("B4").Value = 0.9463
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "0.9096"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9463"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = 0.9276
